$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Cells.Item(11, 1).Value = 9315.02
$ws.Cells.Item(11, 2).Value = 9320.61
$ws.Cells.Item(11, 3).Value = 283.47000000000003
$ws.Cells.Item(11, 4).Value = 283.29000000000002
$ws.Cells.Item(11, 5).Value = $false
$ws.Cells.Item(11, 6).Value = -0.06
$ws.Range("G10").Copy()
$ws.Range("G11").PasteSpecial(-4122)
$ws.Cells.Item(11, 7).Value = 42613.765601851854
$ws.Cells.Item(11, 8).Value = $false

# Row 12
$ws.Cells.Item(12, 1).Value = 9265.65
$ws.Cells.Item(12, 2).Value = 9315.02
$ws.Cells.Item(12, 3).Value = 282.39
$ws.Cells.Item(12, 4).Value = 280.89
$ws.Cells.Item(12, 5).Value = $false
$ws.Cells.Item(12, 6).Value = -0.53
$ws.Range("G10").Copy()
$ws.Range("G12").PasteSpecial(-4122)
$ws.Cells.Item(12, 7).Value = 42614.673043981478
$ws.Cells.Item(12, 8).Value = $false

# Row 13
$ws.Cells.Item(13, 1).Value = 9292.52
$ws.Cells.Item(13, 2).Value = 9265.65
$ws.Cells.Item(13, 3).Value = 280.62
$ws.Cells.Item(13, 4).Value = 281.44
$ws.Cells.Item(13, 5).Value = $false
$ws.Cells.Item(13, 6).Value = 0.28999999999999998
$ws.Range("G10").Copy()
$ws.Range("G13").PasteSpecial(-4122)
$ws.Cells.Item(13, 7).Value = 42615.750219907408
$ws.Cells.Item(13, 8).Value = $true

$excel.CutCopyMode = $false
